# Automatic update of files.
# Bump the "Förändrad" (Changed) date column (C2:C18) forward by one day,
# i.e. from 2023-10-06 (serial 45205) to 2023-10-07 (serial 45206).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C18")
foreach ($cell in $rng.Cells) {
    $cell.Value2 = $cell.Value2 + 1
}
